$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 2 (pushes the existing rows 2-17 down to
# 3-18, preserving each row's own per-cell formatting as it moves).
$ws.Rows.Item(2).Insert(-4121, 1)

# The freshly inserted row 2 inherits the bold/centered header style from
# row 1. Re-stamp it with the plain "data row" formatting used by the rest
# of the table (copied from row 3, which is the row that used to be row 2).
$ws.Range("A3:C3").Copy()
$ws.Range("A2:C2").PasteSpecial(-4122)

# Populate the new row with the camp "return to planet" localization entry.
$ws.Range("A2").Value = "CAMP_RETURN_TO_PLANETE"
$ws.Range("B2").Value = "Return to planet"
$ws.Range("C2").Value = "Retour à la planète"

# Match the author's final selection.
$ws.Range("B3").Select() | Out-Null
